# IoT_Concept_Basic_Architecture.pptx
# Commit: "Changed Windows server to linux for PLC"
#
# 1. "Bare Metall / Hardwareplattform" -> "Any bare metall / Hardwareplattform"
# 2. "SoftPLC VM (Windows Server)" box grows slightly and its caption grows
#    from "(Windows Server)" to
#    "(headless Linux, if no such solution avaliable Windows Server)"
# 3. Neighbouring "Security VM (headless Linux)" box is nudged/resized to
#    make room for the now-larger "SoftPLC VM" box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) "Bare Metall / Hardwareplattform" -> "Any bare metall / Hardwareplattform"
# ---------------------------------------------------------------------------
$bareMetal = $s.Shapes.Item("Rechteck 24")
$bmPara = $bareMetal.TextFrame.TextRange.Paragraphs(1, 1)
$bmPara.Text = "Any bare metall / Hardwareplattform"

# Re-assign each logical chunk back onto itself so the paragraph keeps three
# distinct runs (matching the run-per-word structure of the target deck):
#   "Any bare " | "metall" | " / Hardwareplattform"
$c = $bmPara.Characters(1, 9);  $c.Text = $c.Text   # "Any bare "
$c = $bmPara.Characters(10, 6); $c.Text = $c.Text   # "metall"
$c = $bmPara.Characters(16, 20); $c.Text = $c.Text  # " / Hardwareplattform"

# ---------------------------------------------------------------------------
# 2) SoftPLC VM box: resize + update the "(Windows Server)" caption
# ---------------------------------------------------------------------------
$softPlcBox = $s.Shapes.Item("Rechteck 29")

# Resize/reposition (EMU targets 7058305/1897966/1566878/2143268 converted to
# points; +0.5/12700 nudges past the engine's truncation when dividing back
# to EMU so the round-trip lands exactly on the target EMU value).
$softPlcBox.Left = (7058305 + 0.5) / 12700
$softPlcBox.Top = (1897966 + 0.5) / 12700
$softPlcBox.Width = (1566878 + 0.5) / 12700
$softPlcBox.Height = (2143268 + 0.5) / 12700

# The caption lives in the 4th paragraph of this shape's text body:
#   "Contains only rudimentary ... shutdown."
#   ""
#   "SoftPLC VM"
#   "(Windows Server)"                     <- this one changes
$capPara = $softPlcBox.TextFrame.TextRange.Paragraphs(4, 1)
$capPara.Text = "(headless Linux, if no such solution avaliable Windows Server)"

# Split back into per-word runs, same technique as above, to mirror the
# target deck's run-per-word/err-spellflag structure as closely as COM allows:
#   "(" | "headless" | " Linux, " | "if" | " " | "no" | " such " |
#   "solution" | " " | "avaliable" | " Windows Server)"
$segments = @(
  @(1, 1),    # "("
  @(2, 8),    # "headless"
  @(10, 8),   # " Linux, "
  @(18, 2),   # "if"
  @(20, 1),   # " "
  @(21, 2),   # "no"
  @(23, 6),   # " such "
  @(29, 8),   # "solution"
  @(37, 1),   # " "
  @(38, 9),   # "avaliable"
  @(47, 16)   # " Windows Server)"
)
foreach ($seg in $segments) {
  $c = $capPara.Characters($seg[0], $seg[1])
  $c.Text = $c.Text
}

# ---------------------------------------------------------------------------
# 3) Neighbouring "Security VM" box: resize only, text unchanged
# ---------------------------------------------------------------------------
$securityBox = $s.Shapes.Item("Rechteck 146")
$securityBox.Left = (1736053 + 0.5) / 12700
$securityBox.Top = (845246 + 0.5) / 12700
$securityBox.Width = (1297124 + 0.5) / 12700
$securityBox.Height = (3193332 + 0.5) / 12700
